$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row corrections
$ws.Range("E1").Value = "Title"
$ws.Range("F1").Value = "Country"
$ws.Range("G1").Value = "Phone"
$ws.Range("H1").Value = "Email"

# Row 2 - Aidan McCarron
$ws.Range("A2").Value = "Aidan McCarron"
$ws.Range("C2").Value = "McCarron"
$ws.Range("H2").Value = "aidan.mccarron@mcaleer-rushe.co.uk"

# Row 3 - John Higgins
$ws.Range("A3").Value = "John Higgins "
$ws.Range("C3").Value = "Higgins"
$ws.Range("H3").Value = "john.higgins@mcaleer-rushe.co.uk"

# Row 4 - Michael Yohanis
$ws.Range("A4").Value = "Michael Yohanis "
$ws.Range("C4").Value = "Yohanis"
$ws.Range("H4").Value = "michael.yohanis@mcaleer-rushe.co.uk"

# Row 5 - Declan McLogan
$ws.Range("A5").Value = "Declan McLogan"
$ws.Range("C5").Value = "McLogan"
$ws.Range("H5").Value = "declan.mclogan@mcaleer-rushe.co.uk"

# Row 6 - Lorcan Mulvey
$ws.Range("H6").Value = "lorcan.mulvey@mcaleer-rushe.co.uk"

# Row 7 - Gerald Laverty
$ws.Range("H7").Value = "gerald.laverty@mcaleer-rushe.co.uk"

# Row 9 - Steve Morris
$ws.Range("H9").Value = "steve.morris@mcaleer-rushe.co.uk"

# Row 10 - Lee Robert Gray
$ws.Range("A10").Value = "Lee Robert Gray "
$ws.Range("B10").Value = "Lee"
$ws.Range("C10").Value = "Gray"
$ws.Range("H10").Value = "lee.gray@mcaleer-rushe.co.uk"

# Row 11 - Eamonn Laverty
$ws.Range("C11").Value = "Laverty"
$ws.Range("H11").Value = "eamonn.laverty@thorntonroofing.com"

# Row 12 - Eoin Gormley
$ws.Range("H12").Value = "eoin.gormley@mcaleer-rushe.co.uk"

# Row 13 - Paddy Connolly
$ws.Range("H13").Value = "paddy.connolly@mcaleer-rushe.co.uk"

# Row 14 - Daisy Butterworth
$ws.Range("H14").Value = "daisy.butterworth@mcaleer-rushe.co.uk"

# Row 15 - Sinead Gorman
$ws.Range("C15").Value = "Gorman"
$ws.Range("H15").Value = "sinead.gorman@mcaleer-rushe.co.uk"

# Row 16 - Connor Graham
$ws.Range("H16").Value = "connor.graham@mcaleer-rushe.co.uk"

# Row 17 - Cathal Magee
$ws.Range("H17").Value = "cathal.magee@mcaleer-rushe.co.uk"

# Row 19 - Nina Salandy
$ws.Range("A19").Value = "Nina Salandy "
$ws.Range("C19").Value = "Salandy"
$ws.Range("H19").Value = "nina.salandy@mcaleer-rushe.co.uk"

# Row 20 - Peter Coyle
$ws.Range("H20").Value = "peter.coyle@mcaleer-rushe.co.uk"

# Row 21 - Orran Devine
$ws.Range("H21").Value = "orran.devine@mcaleer-rushe.co.uk"

# Row 22 - Niamh Heneghan
$ws.Range("H22").Value = "niamh.heneghan@mcaleer-rushe.co.uk"
